$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.770.00"
$ws.Range("E2").Value = "  -5.43%  "
$ws.Range("D3").Value = "3.360.97"
$ws.Range("E3").Value = "  -7.33%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "528.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.599"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.01%  "
$ws.Range("D8").Value = "3.353.88"
$ws.Range("E8").Value = "  -7.42%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.618"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.131"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -13.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.27%  "
$ws.Range("D15").Value = "3.913.78"
$ws.Range("E15").Value = "  -6.77%  "
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").Value = "3.375.01"
$ws.Range("E17").Value = "  -6.99%  "
$ws.Range("D18").Value = "64.634.64"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -10.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.957"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -11.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -13.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -17.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "664.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "60.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.95%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("D42").Value = "2.800.84"
$ws.Range("E42").Value = "  -12.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -15.48%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.59%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0617"
$ws.Range("E45").Value = "  -20.47%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0387"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -14.95%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.53%  "
